# Update crypto price/volume figures per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.735.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.723.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4851"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06181"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.727.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06878"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6040"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.461"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9987"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.554.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007116"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.950.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.406"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.571"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.061"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.769"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.365"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07915"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.670"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.595"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9245"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.010"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.591"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3825"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1149"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05374"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.841"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -0.07%  "
